$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sender changes from ECs -> FAPs, values refreshed with new TPM data
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl28"
$ws.Range("C2").Value = "Ccr3"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1395456666666667
$ws.Range("H2").Value = 0.418637
$ws.Range("I2").Value = 0.5708284189068497
$ws.Range("J2").Value = 0.5708284189068498
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3389413333333333
$ws.Range("N2").Value = 1.016824
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.04729779432088888
$ws.Range("R2").Value = 0.425680148888
$ws.Range("S2").Value = 0.5708284189068497
$ws.Range("T2").Value = 0.5708284189068498

# Row 3: sender changes from FAPs -> MuSCs, values refreshed with new TPM data
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Ccl28"
$ws.Range("C3").Value = "Ccr3"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.104916
$ws.Range("H3").Value = 0.314748
$ws.Range("I3").Value = 0.4291715810931503
$ws.Range("J3").Value = 0.4291715810931503
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3389413333333333
$ws.Range("N3").Value = 1.016824
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.035560368928
$ws.Range("R3").Value = 0.320043320352
$ws.Range("S3").Value = 0.4291715810931503
$ws.Range("T3").Value = 0.4291715810931503

# Row 4 (former Resolving-Mac sender row) no longer exists in the updated data
$ws.Rows.Item(4).Delete()
